$p = $ppt.ActivePresentation

# --- Slide 5: update intro text about the "Core" module ---
$slide = $p.Slides.Item(5)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$prefix = $tr.Characters(1, 24)
$prefix.Text = "The Core-Module is for "

# --- Update cached "today" text on every date placeholder (slide master + all layouts) ---
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Datumsplatzhalter*") {
        $sh.TextFrame.TextRange.Text = "01.10.2018"
    }
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Datumsplatzhalter*") {
            $sh.TextFrame.TextRange.Text = "01.10.2018"
        }
    }
}
